$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 23767
$ws.Range("B2").Value = "Paulo Pinto"
$ws.Range("C2").Value = "Operações"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 45078
$ws.Range("G2").Value = 3376.59

# Row 3
$ws.Range("A3").Value = 73215
$ws.Range("B3").Value = "Caio Silva"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45094
$ws.Range("G3").Value = 6376.37

# Row 4
$ws.Range("A4").Value = 69601
$ws.Range("B4").Value = "Henrique Costa"
$ws.Range("C4").Value = "Recursos Humanos"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45098
$ws.Range("G4").Value = 3764.24

# Row 5
$ws.Range("A5").Value = 5428
$ws.Range("B5").Value = "Sr. Luiz Henrique Rezende"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Doença"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45094
$ws.Range("G5").Value = 11529.19

# Row 6
$ws.Range("A6").Value = 27435
$ws.Range("B6").Value = "Benício Nogueira"
$ws.Range("C6").Value = "Operações"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 7186.06

# Row 7
$ws.Range("A7").Value = 27897
$ws.Range("B7").Value = "Ana Nogueira"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45089
$ws.Range("G7").Value = 11211.77

# Row 8
$ws.Range("A8").Value = 15784
$ws.Range("B8").Value = "Marina Nogueira"
$ws.Range("C8").Value = "Atendimento ao Cliente"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 11921.97

# Row 9
$ws.Range("A9").Value = 69625
$ws.Range("B9").Value = "Alana da Rocha"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45100
$ws.Range("G9").Value = 12309.06

# Row 10
$ws.Range("A10").Value = 20697
$ws.Range("B10").Value = "Vitor Gabriel Melo"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45089
$ws.Range("G10").Value = 7230.78

# Row 11
$ws.Range("A11").Value = 19909
$ws.Range("B11").Value = "Marcela Castro"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45090
$ws.Range("G11").Value = 4695.45
